$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Question text" entries for CQ rows 8-18 (column B), previously blank.
$questions = @(
    "What is the specification of an ensemble?",
    "What is the ensemble specification of an ensemble?",
    "What is the region associated to a modeling question?",
    "What are the ensemble specifications of a modeling question?",
    "What are the ensembles of a modeling question?",
    "Which data was selected to run the ensembles in a modeling question?",
    "Which interventions are associated to a modeling question?",
    "Which values in particular were changed in as part of an intervention in the realizations of an ensemble? ",
    "What is the theme of an ensemble?",
    "What are the templates associated to an ensemble?",
    "What are the workflow templates associated to a modeling question?"
)

$startRow = 8
for ($i = 0; $i -lt $questions.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $questions[$i]
}

# Mirror the author's final selection/scroll state left in the saved file.
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("B19").Select()
